# Tarik ha pagato 100 euro il 3/11/2025
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark Tarik's row (row 4) as paid ("✔️") for the months he has now settled:
# Jan 2024 through Apr 2025 (columns N through AC), replacing the old manual
# value (N4) and the projection formulas (O4:AC4) that used to live there.
$ws.Range("N4:AC4").Value = "✔️"

# Reflect the payment in Tarik's running bank balance (row 18).
$ws.Range("B18").Value = 1

# Restore the last-used cell selection as saved by Excel.
[void]$ws.Range("C12").Select()
